# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff" (zh-cn row, and
#   the mirrored Overview/de-de cells that show the same status text).
# - The "Latest Handoff Datetime" timestamp for zh-cn is refreshed.
# - The "Latest HO Xliff Generate Date" timestamp (Overview + de-de, which
#   share the same value) is refreshed.
# - The Status-ish columns are widened so the longer "Ready for handoff"
#   text keeps fitting (mirrors Excel auto-fitting the column after the
#   content grew).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus        = "Ready for handoff"
$newGenerateDate  = "2016-10-20 06:40:41"
$newHandoffDate   = "2016-10-20 06:40:28"

# --- Status: "In Translation" -> "Ready for handoff" -----------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value     = $newStatus
$wsDeDe.Range("C2").Value     = $newStatus

# --- Latest HO Xliff Generate Date ------------------------------------------
$wsOverview.Range("G2").Value = $newGenerateDate
$wsDeDe.Range("H2").Value     = $newGenerateDate

# --- Latest Handoff Datetime (zh-cn only) -----------------------------------
$wsZhCn.Range("H2").Value = $newHandoffDate

# --- Widen the columns that now hold "Ready for handoff" -------------------
# The COM layer quantizes ColumnWidth to sixths of a character; 16.3 is the
# input that lands on the bucket closest to the target rendered width.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.3
